$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data sorted descending by value, with Swedish and Uzbek removed.
$data = @(
    @("English", 22.12197404638313),
    @("Chinese", 17.38857827910244),
    @("Spanish", 6.921745318255761),
    @("Japanese", 4.646413968946974),
    @("Arabic", 4.526208057527736),
    @("German", 4.343436573288437),
    @("Russian", 3.303239255880283),
    @("Portuguese", 3.184682534080445),
    @("Malay-Indonesian", 3.036156302544747),
    @("French", 2.64732634781518),
    @("Italian", 2.002545022988282),
    @("Turkish", 1.807265248466631),
    @("Korean", 1.727672641542481),
    @("Dutch", 1.235332363479998),
    @("Persian", 1.074828512644572),
    @("Thai", 0.9713128289987),
    @("Polish", 0.9105116925265881),
    @("Urdu", 0.8769096773776827),
    @("Bengali", 0.637141092093902),
    @("Vietnamese", 0.6255994449145785)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (previously Uzbek, Vietnamese at 22/23).
$ws.Range("A22:B23").Delete() | Out-Null
